$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LTO7")

$ws.Range("B4").Value = 42.95
$ws.Range("C4").Value = 42.5
$ws.Range("D4").Value = 42.15
$ws.Range("E4").Value = 44.75

$wb.Save()
